$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: date change 01/01/2018 -> 01/01/2021 (keep as plain text, not an actual date,
# so write it via a text formula then paste-values-only to avoid Excel's date auto-detection
# changing the cell's number format/style)
$ws.Range("B8").Formula = "=""01/01/2021"""
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").Formula = "=""01/01/2021"""
$ws.Range("C8").Copy()
$ws.Range("C8").PasteSpecial(-4163)

# Docentes responsáveis: change professor
$ws.Range("B13").Value = "5701460 - Antonio Iacono"
$ws.Range("C13").Value = "5701460 - Antonio Iacono"

# Método:
$ws.Range("B19").Value = "Provas e Trabalhos"
$ws.Range("C19").Value = "Provas e Trabalhos"

# Critério:
$criterio = "M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre projeto de fábricaM = Média de aproveitamento do alunoAprovação com média de aproveitamento maior ou igual a 5,0 e no mínimo 70% de frequência às aulas.A média das provas deve ser maior ou igual a 5,0 (cinco) para que o aluno possa utilizar a nota do Trabalho."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Norma de recuperação:
$norma = "MF = (0,5 M + 0,5 R)M = Média de aproveitamento do aluno, antes da recuperaçãoR = Nota de uma prova de recuperaçãoMF = nota final de aproveitamento, após a recuperaçãoAprovação com média final de aproveitamento maior ou igual a 5,0.A recuperação deverá consistir de uma prova escrita englobando a matéria toda do semestre.Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# Bibliografia:
$bibliografia = "BANZATO, Eduardo et al. Atualidades na armazenagem. São Paulo: IMAM, 2003.BARNES, R.M. Estudo de Movimentos de Tempos: projeto e medida do trabalho. São Paulo, Edgar Blücher, 1977.GURGEL, F.A.C. Administração de recursos materiais e patrimoniais. 2a. Edição. São Paulo. Editora Cengage. 2013. FRANCISCHINI, P.G.; VALLE, C.E. Implantação de Indústrias. Rio de Janeiro, LTC Editora, 1975.LEE, Q et al. Projeto de Instalações e Locais de Trabalho. São Paulo: IMAM, 1998.MOURA, Reinaldo Aparecido. Sistemas e técnicas de movimentação e armazenagem de materiais. IMAM, 2012.NEWMANN, C.; SCALICE, R.K. Projeto de Fábrica e Layout. Rio de Janeiro, Elsevier, 2015.Müther, R. Planejamento do Layout: Sistema SLP. São Paulo, Edgard Blücher, 1978. SLACK, Nigel et al. Administração da produção. São Paulo: Atlas, 8ª ed. 2018.TOMPKINS, James A. et al. Planejamento de instalações. Editora LTC:, 2013."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
